$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "PAP" (Port-au-Prince, Haiti) row, which is row 197.
# This shifts all subsequent rows up by one, shrinking the used range
# from A1:G303 to A1:G302.
$ws.Rows.Item(197).Delete()
